# Estado de Cuenta update: replace the single prior account-statement row with
# the new set of workers/periods, update the summary totals, and move the
# footer signature rows down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Summary header values
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 726623      # VALOR MORA (total)
$ws.Range("C13").Value = 5           # Cant. Trabajadores
$ws.Range("F13").Value = 15          # Cant. Periodos

# ---------------------------------------------------------------------------
# 2) Insert 4 fresh rows above the old "last" (bold) row 29, copying the
#    formatting of row 28 (a normal, non-bold data row) into each of them.
#    This pushes the old row 29 (bold totals-style row) down to row 33 and
#    the footer rows (old 34/35) down to 38/39 automatically.
# ---------------------------------------------------------------------------
$ws.Rows("29:32").Insert()

$styleSrc = $ws.Range("B28:J28")
for ($r = 29; $r -le 32; $r++) {
    $dst = $ws.Range("B" + $r + ":J" + $r)
    $styleSrc.Copy($dst)
}

# ---------------------------------------------------------------------------
# 3) Row 16 now becomes a brand-new worker (KELLY) instead of FABIANA's
#    first period. Keep style (already copied from original file), just
#    change the contents.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "45563382"
$ws.Range("D16").Value = "KELLY DANID OSORIO CASTELLAR"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 60000
$ws.Range("G16").Value = 1500000

# ---------------------------------------------------------------------------
# 4) Rows 17-30: FABIANA ISABEL CASTILLO MENDOZA, one row per overdue period,
#    most recent (2012) first, descending down to 1911. (14 periods total:
#    rows 17-30.) C/D are unchanged from the source file (already FABIANA's
#    doc/name) -- only refresh them defensively, and set E to the period.
# ---------------------------------------------------------------------------
$periods = @("2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","1912","1911")
$r = 17
foreach ($p in $periods) {
    $ws.Range("C" + $r).Value = "1238340678"
    $ws.Range("D" + $r).Value = "FABIANA ISABEL CASTILLO MENDOZA"
    $ws.Range("E" + $r).Value = $p
    $ws.Range("F" + $r).Value = 33125
    $ws.Range("G" + $r).Value = 828116
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 5) Rows 31-33: three more new workers, one overdue period (2507) each.
#    Row 33 keeps the bold/highlighted "last row" style inherited from the
#    old row 29.
# ---------------------------------------------------------------------------
$ws.Range("C31").Value = "1143338851"
$ws.Range("D31").Value = "UBALDO AMETH BANQUEZ GALVIS"
$ws.Range("E31").Value = "2507"
$ws.Range("F31").Value = 76296
$ws.Range("G31").Value = 1907408

$ws.Range("C32").Value = "1047444358"
$ws.Range("D32").Value = "DIANA PAOLA GONZALEZ CUADRADO"
$ws.Range("E32").Value = "2507"
$ws.Range("F32").Value = 69637
$ws.Range("G32").Value = 1740915

$ws.Range("C33").Value = "1047490959"
$ws.Range("D33").Value = "HONEYWELL JOSE SARAVIA SOLANO"
$ws.Range("E33").Value = "2507"
$ws.Range("F33").Value = 56940
$ws.Range("G33").Value = 1423500

# ---------------------------------------------------------------------------
# 6) Recompute the sheet dimension / used range so downstream readers see
#    the right extent (B2:J39).
# ---------------------------------------------------------------------------
$ws.Columns("B:J").AutoFit()
